# Add a new "2022" column (S) to the consumer price index table and
# refresh the 2020 / 2021 figures (columns Q / R) for every region row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- updated values for existing columns Q (2020) and R (2021) ---------
$updates = @{
    5  = @{ Q = 117.60684979252385; R = 113.34848864817617 }
    6  = @{ Q = 114.77319768114526; R = 115.06069350712495 }
    7  = @{ Q = 116.40044011407315; R = 114.29658549692938 }
    8  = @{ Q = 117.53828537152096; R = 113.75761785228545 }
    9  = @{ Q = 117.42206669681742; R = 113.98264089946031 }
    10 = @{ Q = 113.98326995089161; R = 113.92720567782911 }
    11 = @{ Q = 123.488978736909;   R = 114.17226706705155 }
    12 = @{ Q = 118.12340252754679; R = 114.45153946490467 }
    13 = @{ Q = 118.87059844457349; R = 112.69493421065988 }
    14 = @{ Q = 114.06377070452145; R = 113.95067699644588 }
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = $updates[$row].Q
    $ws.Range("R$row").Value = $updates[$row].R
}

# --- new column S (2022) ------------------------------------------------

# Header cell (year label) - same formatting as the rest of the year row.
$ws.Range("S4").Value = 2022
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)

# New data values for each region.
$newValues = @{
    5  = 115.8
    6  = 115.2
    7  = 115.4
    8  = 111.8
    9  = 116.8
    10 = 108.2
    11 = 111
    12 = 115.8
    13 = 117.9
    14 = 112.4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("S$row").Value = $newValues[$row]
}

# Match formatting of the new S cells to the rest of their row (plain
# style for rows 5-13, bottom-bordered "total row" style for row 14).
$ws.Range("A6").Copy()
$ws.Range("S5:S13").PasteSpecial(-4122)

$ws.Range("Q14").Copy()
$ws.Range("S14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- update the active selection to match the source workbook -----------
$ws.Range("T4").Select()
